# "Item menu and cursor"
#
# 1) VARMAP sheet: add a new variable row (row 13) for a pickable-item
#    cursor/selection value.
# 2) SERVICES sheet: rename the GET_ITEM_LIST service to
#    GET_SCENARIO_ITEM_LIST (it now only returns the scenario's item
#    list) and add a brand new GET_PICKED_ITEM_LIST service (row 24) that
#    returns the items the player has already picked up.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: VARMAP
# ---------------------------------------------------------------------
$varmap = $wb.Worksheets.Item("VARMAP")

$varmap.Cells.Item(13, 1).Value = 11          # A13 N
$varmap.Cells.Item(13, 2).Value = "PICKABLE_ITEM_CHOSEN"   # B13 ENUM_ID
$varmap.Cells.Item(13, 3).Value = "GamePickableItem"       # C13 TYPE
$varmap.Cells.Item(13, 4).Value = 1           # D13 SAFE
$varmap.Cells.Item(13, 5).Value = 0           # E13 ARRAY
$varmap.Cells.Item(13, 6).Value = "GamePickableItem.ITEM_PICK_NONE" # F13 DEFVALUE
$varmap.Cells.Item(13, 7).Value = "N"         # G13 SAVE

$varmap.Cells.Item(13, 9).Value  = "-"        # I13  GameMaster
$varmap.Cells.Item(13, 10).Value = "-"        # J13  InputMaster
$varmap.Cells.Item(13, 11).Value = "R"        # K13  LevelMaster
$varmap.Cells.Item(13, 12).Value = "-"        # L13  GraphicsMaster
$varmap.Cells.Item(13, 13).Value = "W"        # M13  GameMenu
$varmap.Cells.Item(13, 14).Value = "R"        # N13  PlayerMaster
$varmap.Cells.Item(13, 15).Value = "R"        # O13  NPCMaster
$varmap.Cells.Item(13, 16).Value = "R"        # P13  ItemMaster
$varmap.Cells.Item(13, 17).Value = "-"        # Q13  GameEventMaster

# Data validation ranges on VARMAP grew by one row (2:12 -> 2:13).
$varmap.Range("D2:D13").Validation.Delete()
$varmap.Range("D2:D13").Validation.Add(3, 1, 1, '"0,1,2"')
$varmap.Range("D2:D13").Validation.ShowInput = $false

$varmap.Range("G2:G13").Validation.Delete()
$varmap.Range("G2:G13").Validation.Add(3, 1, 1, '"N,Y"')
$varmap.Range("G2:G13").Validation.ShowInput = $false

$varmap.Range("I2:Q13").Validation.Delete()
$varmap.Range("I2:Q13").Validation.Add(3, 1, 1, '"R,R E,W,-"')
$varmap.Range("I2:Q13").Validation.ShowInput = $false

# ---------------------------------------------------------------------
# Sheet 2: SERVICES
# ---------------------------------------------------------------------
$services = $wb.Worksheets.Item("SERVICES")

# Row 23: GET_ITEM_LIST -> GET_SCENARIO_ITEM_LIST (renamed in place).
$services.Cells.Item(23, 2).Value = "GET_SCENARIO_ITEM_LIST"
$services.Cells.Item(23, 3).Value = "GET_SCENARIO_ITEM_LIST_DELEGATE"
$services.Cells.Item(23, 4).Value = "LevelMasterClass.GetScenarioItemListService"

# Row 24 (new): GET_PICKED_ITEM_LIST service.
$services.Cells.Item(24, 1).Value = 22
$services.Cells.Item(24, 2).Value = "GET_PICKED_ITEM_LIST"
$services.Cells.Item(24, 3).Value = "GET_PICKED_ITEM_LIST_DELEGATE"
$services.Cells.Item(24, 4).Value = "ItemMasterClass.GetPickedItemListService"
$services.Cells.Item(24, 6).Value  = "-"   # F24
$services.Cells.Item(24, 7).Value  = "-"   # G24
$services.Cells.Item(24, 8).Value  = "-"   # H24
$services.Cells.Item(24, 9).Value  = "-"   # I24
$services.Cells.Item(24, 10).Value = "X"   # J24
$services.Cells.Item(24, 11).Value = "-"   # K24
$services.Cells.Item(24, 12).Value = "-"   # L24
$services.Cells.Item(24, 13).Value = "W"   # M24
$services.Cells.Item(24, 14).Value = "-"   # N24

# Data validation range on SERVICES grew by one row (2:23 -> 2:24).
$services.Range("F2:N24").Validation.Delete()
$services.Range("F2:N24").Validation.Add(3, 1, 1, '"W,X,-"')
$services.Range("F2:N24").Validation.ShowInput = $false

# Column widths: column B (ENUM_ID) gets its own (wider) width, and
# column C (TYPE) widens slightly to fit the new long service names.
$services.Columns.Item(2).ColumnWidth = 23.3
$services.Columns.Item(3).ColumnWidth = 32.8
